# remove AHC examples from US Core
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46: "US Core QuestionnaireResponse Profile" - clear the "Version Retired" (D) and
# "Notes" (E) cells; it's no longer retired in 6.0.0.
$ws.Range("D46").ClearContents()
$ws.Range("E46").ClearContents()

# Row 63: "US Core Extension Questionnaire URI" - same cleanup.
$ws.Range("D63").ClearContents()
$ws.Range("E63").ClearContents()

# Extensions list: remove "US Core Sex For Clinical Use" row content, shifting the
# remaining extension rows up by one and adding a new extension at the end.
$ws.Range("B65").Value = "US Core Tribal Affiliation Extension"
$ws.Range("B66").Value = "US Core Jurisdiction Extension"
$ws.Range("A66").Style = "Normal"
$ws.Range("B67").Value = "US Core USCDI Requirements Extension"

# Add a bold Verdana note in the Notes column for the new extension row.
$ws.Range("E67").Value = "This extension is only used on US Core Profile StructureDefinition elements"
$ws.Range("E67").Font.Bold = $true
$ws.Range("E67").Font.Name = "Verdana"
$ws.Range("E67").Font.Color = 3355443

# Move the active selection to A33 as in the edited workbook.
$ws.Range("A33").Select()
